$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-10 (columns A-T), reflecting recomputed TPM-based
# NATMI ligand-receptor statistics for Fn1-Tnfrsf11b across ECs/FAPs/MuSCs.
$data = @(
    @("ECs", "Fn1", "Tnfrsf11b", "ECs", 3, 1, 6.206015333333333, 18.618046, 0.0150172404156507, 0.0150172404156507, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.028190957994264, 0.02819095799426401, 0.470467679061111, 4.234209111549999, 0.0004233503937473726, 0.0004233503937473726),
    @("ECs", "Fn1", "Tnfrsf11b", "FAPs", 3, 1, 6.206015333333333, 18.618046, 0.0150172404156507, 0.0150172404156507, 3, 1, 2.095195666666667, 6.285587, 0.7791435378093522, 0.7791435378093522, 13.00281643366689, 117.025347903002, 0.01170058582558367, 0.01170058582558367),
    @("ECs", "Fn1", "Tnfrsf11b", "MuSCs", 3, 1, 6.206015333333333, 18.618046, 0.0150172404156507, 0.0150172404156507, 3, 1, 0.5180969999999999, 1.554291, 0.1926655041963838, 0.1926655041963838, 3.215317926153999, 28.937861335386, 0.002893304196319654, 0.002893304196319654),
    @("FAPs", "Fn1", "Tnfrsf11b", "ECs", 3, 1, 385.0524703333334, 1155.157411, 0.9317452840597572, 0.9317452840597571, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.028190957994264, 0.02819095799426401, 29.19018602185278, 262.711674196675, 0.0262667921642822, 0.0262667921642822),
    @("FAPs", "Fn1", "Tnfrsf11b", "FAPs", 3, 1, 385.0524703333334, 1155.157411, 0.9317452840597572, 0.9317452840597571, 3, 1, 2.095195666666667, 6.285587, 0.7791435378093522, 0.7791435378093522, 806.7602672816953, 7260.842405535257, 0.725963316959499, 0.725963316959499),
    @("FAPs", "Fn1", "Tnfrsf11b", "MuSCs", 3, 1, 385.0524703333334, 1155.157411, 0.9317452840597572, 0.9317452840597571, 3, 1, 0.5180969999999999, 1.554291, 0.1926655041963838, 0.1926655041963838, 199.494529722289, 1795.450767500601, 0.179515174935976, 0.179515174935976),
    @("MuSCs", "Fn1", "Tnfrsf11b", "ECs", 3, 1, 22.00088566666667, 66.002657, 0.05323747552459213, 0.05323747552459213, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.028190957994264, 0.02819095799426401, 1.667850474247222, 15.010654268225, 0.001500815436234435, 0.001500815436234435),
    @("MuSCs", "Fn1", "Tnfrsf11b", "FAPs", 3, 1, 22.00088566666667, 66.002657, 0.05323747552459213, 0.05323747552459213, 3, 1, 2.095195666666667, 6.285587, 0.7791435378093522, 0.7791435378093522, 46.09616031162877, 414.865442804659, 0.04147963502426951, 0.0414796350242695),
    @("MuSCs", "Fn1", "Tnfrsf11b", "MuSCs", 3, 1, 22.00088566666667, 66.002657, 0.05323747552459213, 0.05323747552459213, 3, 1, 0.5180969999999999, 1.554291, 0.1926655041963838, 0.1926655041963838, 11.398592861243, 102.587335751187, 0.01025702506408819, 0.01025702506408819)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
